$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025-06-28 05:31"
$ws.Range("B3").Value = "Роман Тарифы FSTA"
$ws.Range("C3").Value = 600
$ws.Range("D3").Value = 200
$ws.Range("E3").Value = 2000
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 7600
